$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'265.77"
$ws.Range("E2").Value = "'1.59%"
$ws.Range("D3").Value = "'26.72"
$ws.Range("E3").Value = "'-1.48%"
$ws.Range("D4").Value = "'4.708"
$ws.Range("E4").Value = "'0.18%"
$ws.Range("D5").Value = "'0.06082"
$ws.Range("E5").Value = "'-1.68%"
$ws.Range("D6").Value = "'6.693"
$ws.Range("E6").Value = "'0.03%"
$ws.Range("D7").Value = "'0.8499"
$ws.Range("E7").Value = "'0.00%"
$ws.Range("D8").Value = "'0.9056"
$ws.Range("E8").Value = "'-1.08%"
$ws.Range("D9").Value = "'0.1409"
$ws.Range("E9").Value = "'-0.10%"
$ws.Range("D10").Value = "'0.05044"
$ws.Range("E10").Value = "'8.79%"
$ws.Range("D11").Value = "'0.07094"
$ws.Range("E11").Value = "'0.18%"
$ws.Range("D12").Value = "'0.03160"
$ws.Range("E12").Value = "'0.24%"
$ws.Range("D13").Value = "'0.09016"
$ws.Range("E13").Value = "'-0.17%"
$ws.Range("D14").Value = "'0.001537"
$ws.Range("E14").Value = "'-0.24%"
$ws.Range("D15").Value = "'0.0006054"
$ws.Range("E15").Value = "'-2.00%"
$ws.Range("D16").Value = "'0.005946"
$ws.Range("E16").Value = "'-3.71%"
$ws.Range("D17").Value = "'3.458"
$ws.Range("E17").Value = "'-0.05%"
$ws.Range("E18").Value = "'0.05%"
$ws.Range("D19").Value = "'2.278"
$ws.Range("E19").Value = "'3.99%"
$ws.Range("D20").Value = "'0.3088"
$ws.Range("E20").Value = "'0.34%"
$ws.Range("D21").Value = "'0.1300"
$ws.Range("E21").Value = "'0.07%"
$ws.Range("D22").Value = "'4.085"
$ws.Range("E22").Value = "'-0.68%"
$ws.Range("D23").Value = "'0.04240"
$ws.Range("E23").Value = "'0.39%"
$ws.Range("E24").Value = "'-2.97%"
$ws.Range("E26").Value = "'0.09%"
$ws.Range("E27").Value = "'5.09%"
$ws.Range("D40").Value = "'0.03921"
$ws.Range("E40").Value = "'0.40%"
$ws.Range("D41").Value = "'0.1113"
$ws.Range("E41").Value = "'0.09%"
$ws.Range("D42").Value = "'0.004172"
$ws.Range("E42").Value = "'1.62%"
$ws.Range("D43").Value = "'0.002111"
$ws.Range("E43").Value = "'-3.30%"
$ws.Range("D44").Value = "'0.01150"
$ws.Range("E44").Value = "'-17.35%"
$ws.Range("D45").Value = "'0.00005105"
$ws.Range("E45").Value = "'-0.73%"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("E46").Value = "'0.05%"
$ws.Range("E48").Value = "'-0.83%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.05%"
